# Update annotations for Ruilin
# - B9 changes from a text "4" to a genuine numeric 4
# - a new row 10 of annotation data is appended (politeness_score stays textual "3")
# - sheet dimension grows to A1:H10 automatically as cells are written

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 was stored as text "4"; make it a real number.
$ws.Range("B9").Value = 4

# New row 10.
$ws.Range("A10").Value = "Ruilin"

# B10's politeness_score must stay a text value "3" (like similar cells elsewhere
# in the column), not get auto-coerced into a number. Enter it quote-prefixed
# (like typing '3 into Excel) then restore the cell to the default "Normal"
# style so no stray number-format/quote-prefix styling is left behind.
$ws.Range("B10").Value = "'3"
$ws.Range("B10").Style = "Normal"

$ws.Range("C10").Value = "无"
$ws.Range("D10").Value = "QSN"
$ws.Range("E10").Value = "RES"
$ws.Range("F10").Value = "a5228610-fe6d-4383-b598-a7c34c3b8714"
$ws.Range("G10").Value = "HyRnez-RW_annotated.xlsx"
$ws.Range("H10").Value = "Why is this result not compared to in Table 1?"
